# Regenerate save_data to use K (strikeouts-derived) instead of Strike# in column G.
# Column G (header "K") values are rewritten with the newly calculated s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    11 = 2
    12 = 0
    13 = 2
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 2
    19 = 0
    20 = 2
    21 = 2
    22 = 2
    23 = 0
    24 = 1
    25 = 0
    26 = 2
    27 = 2
    28 = 2
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 1
    34 = 2
    35 = 1
    36 = 2
    37 = 0
    38 = 1
    40 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
